$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-13, columns B, C, D, F
$ws.Range("B2").Value = "NSE:ADROITINFO"
$ws.Range("C2").Value = "NSE:AGI"
$ws.Range("D2").Value = "NSE:ABCAPITAL"
$ws.Range("F2").Value = "NSE:OFSS"

$ws.Range("B3").Value = "NSE:AXISILVER"
$ws.Range("C3").Value = "NSE:ALKYLAMINE"
$ws.Range("D3").Value = "NSE:GODREJPROP"
$ws.Range("F3").Value = ""

$ws.Range("B4").Value = "NSE:BAJAJFINSV"
$ws.Range("C4").Value = "NSE:DEEPAKFERT"

$ws.Range("B5").Value = "NSE:BAJFINANCE"
$ws.Range("C5").Value = "NSE:DOLATALGO"

$ws.Range("B6").Value = "NSE:BFUTILITIE"
$ws.Range("C6").Value = "NSE:FINPIPE"

$ws.Range("B7").Value = "NSE:CAMPUS"
$ws.Range("C7").Value = "NSE:HIRECT"

$ws.Range("B8").Value = "NSE:CANFINHOME"
$ws.Range("C8").Value = "NSE:LAXMICOT"

$ws.Range("B9").Value = "NSE:LICNETFSEN"
$ws.Range("C9").Value = "NSE:MOREPENLAB"

$ws.Range("B10").Value = "NSE:OFSS"
$ws.Range("C10").Value = "NSE:NIITLTD"

$ws.Range("B11").Value = "NSE:SAGCEM"
$ws.Range("C11").Value = "NSE:ORIENTHOT"

$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:PAKKA"

$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:PASUPTAC"

# Delete rows 14-22 (the old rows 14 through 22, which contained data that no longer exists)
$ws.Range("A14:F22").EntireRow.Delete()
